$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'64.252.32"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -3.03%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.133.59"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -2.46%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.01%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'607.75"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.17%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'147.31"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -5.00%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.06%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'3.130.94"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -2.52%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.526"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -3.46%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  -5.14%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'5.54"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -2.80%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.474"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -5.15%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.0000256"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -4.00%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'36.39"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -4.78%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'3.646.11"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -2.45%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'64.197.05"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -3.20%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.02%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'3.129.20"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -2.55%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'6.93"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -4.26%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'478.94"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -5.36%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'14.52"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -4.44%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'0.705"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -3.04%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'7.69"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -3.54%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'13.71"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -5.40%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'83.29"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -1.85%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +0.04%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'2.92"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -2.50%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  -5.88%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  -5.28%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'0.122"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -19.12%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'6.78"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -1.83%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'0.999"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -0.20%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  -5.40%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'26.42"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -6.40%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'1.11"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -4.79%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'6.05"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -4.98%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'54.52"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -1.47%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'3.08"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +2.66%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.0₃0726"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -4.99%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'450.88"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -9.48%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.0397"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -4.76%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  -5.10%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'8.39"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -3.57%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'2.861.22"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -2.00%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.270"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -7.92%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'  -7.43%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'26.46"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -5.29%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'0.999"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E49').Value = "'  -2.46%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'2.31"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -2.97%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'118.74"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -2.38%  "
$ws.Range('E51').Style = 'Normal'
